$wb = $excel.ActiveWorkbook

# ----- Sheet 1: quality_comparison -----
$ws1 = $wb.Worksheets.Item(1)

# C1, D1 -> thin top/bottom border (C1), thin top/bottom/right border (D1)
foreach ($addr in @("C1", "D1")) {
    $r = $ws1.Range($addr)
    $r.ClearFormats()
    $r.Borders.Item(7).LineStyle = -4142    # xlEdgeLeft -> none
    $r.Borders.Item(8).LineStyle = 1        # xlEdgeTop -> continuous
    $r.Borders.Item(8).Weight = 2           # xlThin
    $r.Borders.Item(9).LineStyle = 1        # xlEdgeBottom -> continuous
    $r.Borders.Item(9).Weight = 2           # xlThin
    if ($addr -eq "D1") {
        $r.Borders.Item(10).LineStyle = 1   # xlEdgeRight -> continuous
        $r.Borders.Item(10).Weight = 2      # xlThin
    } else {
        $r.Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none
    }
}

$ws1.Range("C2").Value = "approach"

# ----- Sheet 2: computational_comparison -----
$ws2 = $wb.Worksheets.Item(2)

foreach ($addr in @("C1", "D1", "F1", "G1")) {
    $r = $ws2.Range($addr)
    $r.ClearFormats()
    $r.Borders.Item(7).LineStyle = -4142    # xlEdgeLeft -> none
    $r.Borders.Item(8).LineStyle = 1        # xlEdgeTop -> continuous
    $r.Borders.Item(8).Weight = 2           # xlThin
    $r.Borders.Item(9).LineStyle = 1        # xlEdgeBottom -> continuous
    $r.Borders.Item(9).Weight = 2           # xlThin
    if ($addr -eq "D1" -or $addr -eq "G1") {
        $r.Borders.Item(10).LineStyle = 1   # xlEdgeRight -> continuous
        $r.Borders.Item(10).Weight = 2      # xlThin
    } else {
        $r.Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none
    }
}

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

$ws2.Range("G5").ClearContents()
